$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.437.79'
$ws.Range('E2').Value = '  +1.40%  '
$ws.Range('D3').Value = '1.667.56'
$ws.Range('E3').Value = '  +1.34%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.0000'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.08'
$ws.Range('E5').Value = '  +1.80%  '
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3963'
$ws.Range('E7').Value = '  +1.40%  '
$ws.Range('E8').Value = '  +2.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '52.03'
$ws.Range('E9').Value = '  +6.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.395'
$ws.Range('E10').Value = '  +3.42%  '
$ws.Range('E11').Value = '  +0.08%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08580'
$ws.Range('E12').Value = '  +1.33%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '24.55'
$ws.Range('E13').Value = '  +2.96%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.322'
$ws.Range('E14').Value = '  +2.73%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.949'
$ws.Range('E15').Value = '  +6.14%  '
$ws.Range('E16').Value = '  +4.35%  '
$ws.Range('D17').Value = '1.657.08'
$ws.Range('E17').Value = '  +0.77%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '95.25'
$ws.Range('E18').Value = '  +0.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.07009'
$ws.Range('E19').Value = '  +0.56%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '20.60'
$ws.Range('E20').Value = '  -0.62%  '
$ws.Range('E21').Value = '  +0.99%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9984'
$ws.Range('E22').Value = '  -0.31%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.75'
$ws.Range('E23').Value = '  +0.90%  '
$ws.Range('D24').Value = '24.399.25'
$ws.Range('E24').Value = '  +1.26%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.458'
$ws.Range('E25').Value = '  +5.67%  '
$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.094'
$ws.Range('E26').Value = '  +14.41%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.54'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '157.88'
$ws.Range('E28').Value = '  -0.12%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '142.88'
$ws.Range('E29').Value = '  +0.99%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.443'
$ws.Range('E30').Value = '  +2.77%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.083'
$ws.Range('E31').Value = '  -7.71%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.548'
$ws.Range('E32').Value = '  +3.57%  '
$ws.Range('D33').Value = '1.844.99'
$ws.Range('E33').Value = '  +0.94%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.065'
$ws.Range('E34').Value = '  +11.25%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.03069'
$ws.Range('E35').Value = '  +5.36%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.08267'
$ws.Range('E36').Value = '  +3.07%  '
$ws.Range('E37').Value = '  -0.06%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '11.15'
$ws.Range('E38').Value = '  +11.91%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2767'
$ws.Range('E39').Value = '  +2.57%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.09276'
$ws.Range('E40').Value = '  +0.79%  '
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '13.83'
$ws.Range('E41').Value = '  +5.80%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.7711'
$ws.Range('E42').Value = '  +1.35%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.449'
$ws.Range('E43').Value = '  -0.95%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.57'
$ws.Range('E44').Value = '  +2.88%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.7135'
$ws.Range('E45').Value = '  +3.20%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.542'
$ws.Range('E46').Value = '  +2.42%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.136'
$ws.Range('E47').Value = '  +0.85%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.000'
$ws.Range('E48').Value = '  -0.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.08441'
$ws.Range('E49').Value = '  +1.08%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '136.81'
$ws.Range('E50').Value = '  +2.18%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.269'
$ws.Range('E51').Value = '  +0.59%  '
